# Auto-generated script applying the Atomos_Profits market-data refresh diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value2 = 2215.6
$ws.Cells.Item(86, 9).Value2 = 2217.3333
$ws.Cells.Item(86, 10).Value2 = 2200
$ws.Cells.Item(86, 11).Value2 = 2217.3333
$ws.Cells.Item(86, 12).Value2 = 2200
$ws.Cells.Item(86, 13).Value2 = -1094.3333
$ws.Cells.Item(86, 14).Value2 = -4446

$ws.Cells.Item(89, 8).Value2 = 2215.6
$ws.Cells.Item(89, 9).Value2 = 2217.3333
$ws.Cells.Item(89, 10).Value2 = 2200
$ws.Cells.Item(89, 11).Value2 = 11086.6665
$ws.Cells.Item(89, 12).Value2 = 11000
$ws.Cells.Item(89, 13).Value2 = -5470.666499999999
$ws.Cells.Item(89, 14).Value2 = -22232

$ws.Cells.Item(106, 8).Value2 = 1925.2963
$ws.Cells.Item(106, 9).Value2 = 1292.2
$ws.Cells.Item(106, 10).Value2 = 2716.6667
$ws.Cells.Item(106, 11).Value2 = 1292.2
$ws.Cells.Item(106, 12).Value2 = 2716.6667
$ws.Cells.Item(106, 13).Value2 = -661.2
$ws.Cells.Item(106, 14).Value2 = -3978.6667

$ws.Cells.Item(111, 8).Value2 = 791.38464
$ws.Cells.Item(111, 9).Value2 = 762.5454999999999
$ws.Cells.Item(111, 11).Value2 = 2287.6365
$ws.Cells.Item(111, 13).Value2 = 779.3635000000004

$ws.Cells.Item(116, 8).Value2 = 3968.84
$ws.Cells.Item(116, 9).Value2 = 3145.7693
$ws.Cells.Item(116, 10).Value2 = 4860.5
$ws.Cells.Item(116, 11).Value2 = 3145.7693
$ws.Cells.Item(116, 12).Value2 = 4860.5
$ws.Cells.Item(116, 13).Value2 = 296.2307000000001
$ws.Cells.Item(116, 14).Value2 = -11744.5

$ws.Cells.Item(127, 8).Value2 = 996.8570999999999
$ws.Cells.Item(127, 9).Value2 = 1575
$ws.Cells.Item(127, 10).Value2 = 952.38464
$ws.Cells.Item(127, 11).Value2 = 4725
$ws.Cells.Item(127, 12).Value2 = 2857.15392
$ws.Cells.Item(127, 13).Value2 = 235
$ws.Cells.Item(127, 14).Value2 = -12777.15392

$ws.Cells.Item(129, 8).Value2 = 5001054.5
$ws.Cells.Item(129, 10).Value2 = 998.45654
$ws.Cells.Item(129, 12).Value2 = 2995.36962
$ws.Cells.Item(129, 14).Value2 = -12995.36962

$ws.Cells.Item(134, 8).Value2 = 21870.525
$ws.Cells.Item(134, 10).Value2 = 21870.525
$ws.Cells.Item(134, 12).Value2 = 21870.525
$ws.Cells.Item(134, 14).Value2 = -32010.525

$ws.Cells.Item(135, 8).Value2 = 579.4138
$ws.Cells.Item(135, 9).Value2 = 400.1111
$ws.Cells.Item(135, 11).Value2 = 3600.9999
$ws.Cells.Item(135, 13).Value2 = -1065.9999

$ws.Cells.Item(136, 8).Value2 = 29400
$ws.Cells.Item(136, 10).Value2 = 29400
$ws.Cells.Item(136, 12).Value2 = 29400
$ws.Cells.Item(136, 14).Value2 = -39600

$ws.Cells.Item(138, 8).Value2 = 4833.663
$ws.Cells.Item(138, 9).Value2 = 2567.6843
$ws.Cells.Item(138, 10).Value2 = 5423.4385
$ws.Cells.Item(138, 11).Value2 = 7703.0529
$ws.Cells.Item(138, 12).Value2 = 16270.3155
$ws.Cells.Item(138, 13).Value2 = -2563.0529
$ws.Cells.Item(138, 14).Value2 = -26550.3155

$ws.Cells.Item(139, 8).Value2 = 28909.092
$ws.Cells.Item(139, 10).Value2 = 28909.092
$ws.Cells.Item(139, 12).Value2 = 28909.092
$ws.Cells.Item(139, 14).Value2 = -39189.092

$ws.Cells.Item(140, 8).Value2 = 0
$ws.Cells.Item(140, 10).Value2 = 0
$ws.Cells.Item(140, 12).Value2 = 0
$ws.Cells.Item(140, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value2 = 3959.7368
$ws.Cells.Item(61, 9).Value2 = 2514.9167
$ws.Cells.Item(61, 10).Value2 = 6436.5713
$ws.Cells.Item(61, 11).Value2 = 2514.9167
$ws.Cells.Item(61, 12).Value2 = 6436.5713
$ws.Cells.Item(61, 13).Value2 = -2302.9167
$ws.Cells.Item(61, 14).Value2 = -6860.5713

$ws.Cells.Item(136, 8).Value2 = 3959.7368
$ws.Cells.Item(136, 9).Value2 = 2514.9167
$ws.Cells.Item(136, 10).Value2 = 6436.5713
$ws.Cells.Item(136, 11).Value2 = 7544.750100000001
$ws.Cells.Item(136, 12).Value2 = 19309.7139
$ws.Cells.Item(136, 13).Value2 = -4994.750100000001
$ws.Cells.Item(136, 14).Value2 = -24409.7139

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value2 = 2785.012
$ws.Cells.Item(31, 9).Value2 = 1844.9807
$ws.Cells.Item(31, 10).Value2 = 4312.5625
$ws.Cells.Item(31, 11).Value2 = 1844.9807
$ws.Cells.Item(31, 12).Value2 = 4312.5625
$ws.Cells.Item(31, 13).Value2 = -1549.9807
$ws.Cells.Item(31, 14).Value2 = -4902.5625

$ws.Cells.Item(34, 8).Value2 = 2785.012
$ws.Cells.Item(34, 9).Value2 = 1844.9807
$ws.Cells.Item(34, 10).Value2 = 4312.5625
$ws.Cells.Item(34, 11).Value2 = 1844.9807
$ws.Cells.Item(34, 12).Value2 = 4312.5625
$ws.Cells.Item(34, 13).Value2 = -1642.9807
$ws.Cells.Item(34, 14).Value2 = -4716.5625

$ws.Cells.Item(132, 8).Value2 = 2509.7646
$ws.Cells.Item(132, 9).Value2 = 2187.1853
$ws.Cells.Item(132, 10).Value2 = 3754
$ws.Cells.Item(132, 11).Value2 = 6561.5559
$ws.Cells.Item(132, 12).Value2 = 11262
$ws.Cells.Item(132, 13).Value2 = -4031.5559
$ws.Cells.Item(132, 14).Value2 = -16322

$ws.Cells.Item(134, 8).Value2 = 16669982
$ws.Cells.Item(134, 9).Value2 = 27780988
$ws.Cells.Item(134, 10).Value2 = 3474.4167
$ws.Cells.Item(134, 11).Value2 = 83342964
$ws.Cells.Item(134, 12).Value2 = 10423.2501
$ws.Cells.Item(134, 13).Value2 = -83340429
$ws.Cells.Item(134, 14).Value2 = -15493.2501

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value2 = 494.1
$ws.Cells.Item(4, 10).Value2 = 900.8
$ws.Cells.Item(4, 12).Value2 = 2702.4
$ws.Cells.Item(4, 14).Value2 = -2926.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value2 = 40272.332
$ws.Cells.Item(102, 9).Value2 = 2450.6667
$ws.Cells.Item(102, 10).Value2 = 172648.17
$ws.Cells.Item(102, 11).Value2 = 2450.6667
$ws.Cells.Item(102, 12).Value2 = 172648.17
$ws.Cells.Item(102, 13).Value2 = -828.6667000000002
$ws.Cells.Item(102, 14).Value2 = -175892.17

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value2 = 1255.1428
$ws.Cells.Item(22, 9).Value2 = 517.5
$ws.Cells.Item(22, 10).Value2 = 1550.2
$ws.Cells.Item(22, 11).Value2 = 517.5
$ws.Cells.Item(22, 12).Value2 = 1550.2
$ws.Cells.Item(22, 13).Value2 = -222.5
$ws.Cells.Item(22, 14).Value2 = -2140.2

$ws.Cells.Item(27, 8).Value2 = 1255.1428
$ws.Cells.Item(27, 9).Value2 = 517.5
$ws.Cells.Item(27, 10).Value2 = 1550.2
$ws.Cells.Item(27, 11).Value2 = 517.5
$ws.Cells.Item(27, 12).Value2 = 1550.2
$ws.Cells.Item(27, 13).Value2 = -410.5
$ws.Cells.Item(27, 14).Value2 = -1764.2

$ws.Cells.Item(40, 8).Value2 = 5923.2354
$ws.Cells.Item(40, 9).Value2 = 7509
$ws.Cells.Item(40, 10).Value2 = 3657.8572
$ws.Cells.Item(40, 11).Value2 = 7509
$ws.Cells.Item(40, 12).Value2 = 3657.8572
$ws.Cells.Item(40, 13).Value2 = -7373
$ws.Cells.Item(40, 14).Value2 = -3929.8572

$ws.Cells.Item(82, 8).Value2 = 4286.143
$ws.Cells.Item(82, 9).Value2 = 1750
$ws.Cells.Item(82, 10).Value2 = 5300.6
$ws.Cells.Item(82, 11).Value2 = 1750
$ws.Cells.Item(82, 12).Value2 = 5300.6
$ws.Cells.Item(82, 13).Value2 = -1389
$ws.Cells.Item(82, 14).Value2 = -6022.6

$ws.Cells.Item(85, 8).Value2 = 4286.143
$ws.Cells.Item(85, 9).Value2 = 1750
$ws.Cells.Item(85, 10).Value2 = 5300.6
$ws.Cells.Item(85, 11).Value2 = 1750
$ws.Cells.Item(85, 12).Value2 = 5300.6
$ws.Cells.Item(85, 13).Value2 = -502
$ws.Cells.Item(85, 14).Value2 = -7796.6

$ws.Cells.Item(132, 8).Value2 = 3014.8823
$ws.Cells.Item(132, 9).Value2 = 2323.9412
$ws.Cells.Item(132, 10).Value2 = 3705.8235
$ws.Cells.Item(132, 11).Value2 = 6971.823600000001
$ws.Cells.Item(132, 12).Value2 = 11117.4705
$ws.Cells.Item(132, 13).Value2 = -4441.823600000001
$ws.Cells.Item(132, 14).Value2 = -16177.4705

$ws.Cells.Item(136, 8).Value2 = 3461.7693
$ws.Cells.Item(136, 9).Value2 = 2535.8572
$ws.Cells.Item(136, 10).Value2 = 5818.636
$ws.Cells.Item(136, 11).Value2 = 7607.571599999999
$ws.Cells.Item(136, 12).Value2 = 17455.908
$ws.Cells.Item(136, 13).Value2 = -5057.571599999999
$ws.Cells.Item(136, 14).Value2 = -22555.908

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value2 = 2232
$ws.Cells.Item(81, 9).Value2 = 2062.5
$ws.Cells.Item(81, 10).Value2 = 2425.7144
$ws.Cells.Item(81, 11).Value2 = 4125
$ws.Cells.Item(81, 12).Value2 = 4851.4288
$ws.Cells.Item(81, 13).Value2 = -3064
$ws.Cells.Item(81, 14).Value2 = -6973.4288

$ws.Cells.Item(84, 8).Value2 = 2232
$ws.Cells.Item(84, 9).Value2 = 2062.5
$ws.Cells.Item(84, 10).Value2 = 2425.7144
$ws.Cells.Item(84, 11).Value2 = 20625
$ws.Cells.Item(84, 12).Value2 = 24257.144
$ws.Cells.Item(84, 13).Value2 = -15321
$ws.Cells.Item(84, 14).Value2 = -34865.144

$ws.Cells.Item(107, 8).Value2 = 1330.2222
$ws.Cells.Item(107, 9).Value2 = 246.5
$ws.Cells.Item(107, 10).Value2 = 10000
$ws.Cells.Item(107, 11).Value2 = 739.5
$ws.Cells.Item(107, 12).Value2 = 30000
$ws.Cells.Item(107, 13).Value2 = 1180.5
$ws.Cells.Item(107, 14).Value2 = -33840
